$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain text in this sheet (values like
# "43.821.19" are not valid numbers). Force Text format on the specific cells
# we touch so Excel does not reinterpret numeric-looking strings as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.837.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.349.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.01"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.46"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.904"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.346.79"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.826.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.54"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("B23").Value = "WEMIXToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.86"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.63"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.16"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.128"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0741"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.76"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.44"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.57"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +18.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0273"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +13.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.63"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.05"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.106"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.198"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.23"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.83%  "
